# Update web reports - 2025-12-31 16:42:34
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weekly hours")

# --- Header row additions: new "Engineer/ supervisor" column (X) ---
$ws.Range("X1").Value = "Engineer/ supervisor  "

# --- Row 3 (names): Brian M / Vu T swap + new Tony M column ---
$ws.Range("O3").Value = "Brian M"
$ws.Range("X3").Value = "Vu T"
# Copy X3's header formatting onto Y3 before stamping the new name so the
# style (s="9") matches the rest of the header row, then overwrite the value.
$ws.Range("X3").Copy($ws.Range("Y3"))
$ws.Range("Y3").Value = "Tony M "

# --- Week 32 (row 35): move Vu T's hours from column O to the new column X ---
$ws.Range("O35").Value = ""
$ws.Range("X35").Value = 32

# --- Week 33 (row 36): move Vu T's hours from column O to the new column X ---
$ws.Range("O36").Value = ""
$ws.Range("X36").Value = 40

# --- Week 34 (row 37): Luis Z (I) logged 1 hour ---
$ws.Range("I37").Value = 1

# --- Week 35 (row 38): additional hours logged ---
$ws.Range("I38").Value = 1
$ws.Range("P38").Value = 40
$ws.Range("Q38").Value = 40
$ws.Range("R38").Value = 40
$ws.Range("W38").Value = 40

# --- Week 36 (row 39): additional hours logged ---
$ws.Range("L39").Value = 40
$ws.Range("P39").Value = 40
$ws.Range("Q39").Value = 40
$ws.Range("R39").Value = 40
$ws.Range("V39").Value = 40
$ws.Range("W39").Value = 40

# --- Week 37 (row 40): additional hours logged + new Tony M column ---
$ws.Range("L40").Value = 40
$ws.Range("P40").Value = 40
$ws.Range("Q40").Value = 40
$ws.Range("R40").Value = 40
$ws.Range("V40").Value = 40
$ws.Range("W40").Value = 40
$ws.Range("Y40").Value = 40
